$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.382.77'
$ws.Range("E2").Value = '  -0.03%  '
$ws.Range("D3").Value = '1.848.07'
$ws.Range("E3").Value = '  -0.01%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9993'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.23'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.10%  '
$ws.Range("E6").Value = '  -0.21%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.0000'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07635'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.19%  '
$ws.Range("E9").Value = '  -1.13%  '
$ws.Range("E10").Value = '  +0.88%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07740'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.04%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.032'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.53%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.6787'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.06%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.00001060'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.35%  '
$ws.Range("E15").Value = '  -0.40%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.154'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.04%  '
$ws.Range("D17").Value = '29.402.82'
$ws.Range("E17").Value = '  -0.04%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '227.51'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.51%  '
$ws.Range("E19").Value = '  -0.75%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9997'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.03%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.494'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.41%  '
$ws.Range("E22").Value = '  -0.03%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '158.29'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.33%  '
$ws.Range("E24").Value = '  -0.35%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.402'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.26%  '
$ws.Range("E26").Value = '  +0.31%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.376'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +5.16%  '
$ws.Range("E28").Value = '  -0.59%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.05600'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.80%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.116'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.10%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.076'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.98%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.836'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.86%  '
$ws.Range("E33").Value = '  +0.47%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.6939'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.23%  '
$ws.Range("E35").Value = '  -0.14%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.01803'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.32%  '
$ws.Range("D37").Value = '1.229.13'
$ws.Range("E37").Value = '  -0.32%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.710'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.27%  '
$ws.Range("E39").Value = '  -1.16%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9048'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.28%  '
$ws.Range("E41").Value = '  +0.04%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '101.52'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.02%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '66.06'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.08%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '7.173'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.20%  '
$ws.Range("B45").Value = 'BabyDogeCoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00000000119'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.82%  '
$ws.Range("B46").Value = 'TheSandbox'
$ws.Range("C46").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4010'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.12%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.002'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.21%  '
$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.673'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.78%  '
$ws.Range("B49").Value = 'Algorand'
$ws.Range("C49").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.1138'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.27%  '
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05702'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.14%  '
$ws.Range("B51").Value = 'Mantle'
$ws.Range("C51").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4627'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.07%  '
